{"js": "// The template's mustache-style placeholders (\"{{ ... }}\") and several\n// runs of plain prose had been split across multiple Word runs (often to\n// bracket a spelling/grammar \"wavy underline\" flagged by the proofer).\n// The authoritative edit re-saves the document so those runs are no\n// longer artificially split - the only *visible* text change is that the\n// blank-signature line grows from 6 underscores to 9 underscores, and the\n// stray tab that used to separate \"on: \" from \"___ Prosecutor's Office\"\n// is dropped (the extra underscores now occupy that space).\n\nconst body = context.document.body;\n\n// 1) Grow \"Dep. Clerk ______ on:\" -> \"Dep. Clerk _________ on:\" (+3 \"_\").\nconst sigResults = body.search(\"Dep. Clerk ______ on:\", { matchCase: true });\nsigResults.load(\"text\");\nawait context.sync();\n\nif (sigResults.items.length > 0) {\n  sigResults.items[0].insertText(\"Dep. Clerk _________ on:\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the tab character that used to sit between \"on: \" and the\n//    \"___ Prosecutor's Office\" line - the paragraph now flows directly\n//    from one piece of text into the next.\nconst tabResults = body.search(\"on: \\t___ Prosecutor\", { matchCase: true });\ntabResults.load(\"text\");\nawait context.sync();\n\nif (tabResults.items.length > 0) {\n  tabResults.items[0].insertText(\"on: ___ Prosecutor\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The template's mustache-style placeholders (\"{{ ... }}\") and several\n# runs of plain prose had been split across multiple Word runs (often to\n# bracket a spelling/grammar \"wavy underline\" flagged by the proofer).\n# The authoritative edit re-saves the document so those runs are no\n# longer artificially split - the only *visible* text change is that the\n# blank signature line grows from 6 underscores to 9 underscores, and the\n# stray tab that used to separate \"on: \" from \"___ Prosecutor's Office\"\n# is dropped (the extra underscores now occupy that space).\n\n$d = $word.ActiveDocument\n\n# 1) Grow \"Dep. Clerk ______ on:\" -> \"Dep. Clerk _________ on:\" (+3 \"_\").\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n    \"Dep. Clerk ______ on:\",  # FindText\n    $false,                   # MatchCase\n    $false,                   # MatchWholeWord\n    $false,                   # MatchWildcards\n    $false,                   # MatchSoundsLike\n    $false,                   # MatchAllWordForms\n    $true,                    # Forward\n    1,                        # Wrap (wdFindContinue)\n    $false,                   # Format\n    \"Dep. Clerk _________ on:\", # ReplaceWith\n    2                         # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) Remove the tab character that used to sit between \"on: \" and the\n#    \"___ Prosecutor's Office\" line - the paragraph now flows directly\n#    from one piece of text into the next.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    \"on: ^t___ Prosecutor\",   # FindText (\"^t\" = tab char)\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"on: ___ Prosecutor\",\n    2\n) | Out-Null\n"}
